# Add results for the additional "Unpadded" experiment.
# This mirrors the existing CLS/MEAN/MEAN_NO_CLS (D/E columns) result blocks
# for the Audio/SpeechCommands and Text/IMDB sections into new F/G/H columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Audio / SpeechCommands / MEAN block (rows 36-38) ---
$ws.Range("F36").Value = "Unpadded"
$ws.Range("G37").Value = "Top1"
$ws.Range("H37").Value = "0.8236"
$ws.Range("G38").Value = "Top5"
$ws.Range("H38").Value = "0.8962"

# --- Audio / SpeechCommands / MEAN_NO_CLS block (rows 40-42) ---
$ws.Range("F40").Value = "Unpadded"
$ws.Range("G41").Value = "Top1"
$ws.Range("H41").Value = "0.8240"
$ws.Range("G42").Value = "Top5"
$ws.Range("H42").Value = "0.8972"

# --- Text / IMDB / MEAN block (rows 50-51) ---
$ws.Range("F50").Value = "Unpadded"
$ws.Range("G51").Value = "Top1"
$ws.Range("H51").Value = "0.1488"

# --- Text / IMDB / MEAN_NO_CLS block (rows 54-55) ---
$ws.Range("F54").Value = "Unpadded"
$ws.Range("G55").Value = "Top1"
$ws.Range("H55").Value = "0.1484"

# --- Update view: scroll + selection to match the saved state ---
$ws.Range("A25").Select()
$sheetView = $ws.Application
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("I30").Select()
